# Update column G ("K" - strikeouts) values to the newly regenerated
# save_data values (using actual K instead of the old Strike# figure).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 6
    3  = 6
    4  = 7
    5  = 4
    6  = 6
    7  = 1
    8  = 3
    9  = 4
    10 = 7
    11 = 9
    12 = 9
    13 = 10
    14 = 8
    15 = 5
    16 = 7
    17 = 8
    18 = 4
    19 = 5
    20 = 3
    21 = 6
    22 = 5
    23 = 3
    24 = 0
    25 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
